# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Splits the old citation block (A45 text + A46 hyperlinked URL + A51 long
# combined citation) into separate plain-text lines and drops the
# hyperlink:
#   A45: blank separator line (new)
#   A46: "Bank of Jamaica, Financing ..." (was A45, unchanged style)
#   A47: blank separator line (was A46's slot, now plain/blank)
#   A48: the bare URL as plain text (was A46's hyperlinked value, hyperlink removed)
#   A51: "Bank of Jamaica" (unchanged, was A50)
#   A52: "Bank of Jamaica" (was A51's long citation, now shortened)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = 'http://www.boj.org.jm/uploads/pdf/papers_pamphlets/papers_pamphlets_Financiang_the_MSME_Sector_in_Jamaica_-_Constraints_and_Prospects_for_Leasing,_Factoring_and_Microfinance.pdf'
$bankName = 'Bank of Jamaica'

# Remove the existing hyperlink (anchored on A46) before we move things
# around, so it doesn't tag along to the wrong cell.
$ws.Hyperlinks.Delete()

# Insert a new blank row above the old A45 - pushes the old A45..A51
# block down by one row (A45->A46, A46->A47, A47->A48, A50->A51, A51->A52).
$ws.Rows(45).Insert()

# The row that used to hold the hyperlinked URL (now at row 47) becomes a
# blank separator line, and the row that used to be blank (now at row 48)
# picks up the bare URL as plain text.
$ws.Range("A47").Value = ""
$ws.Range("A48").Value = $url

# The long combined citation (now at row 52) is replaced with the short
# "Bank of Jamaica" line.
$ws.Range("A52").Value = $bankName
